# The mod's "Faceless" boss/track is being renamed to "Mirthless" (part of
# the ongoing "removing sculkwood" rework). Update the Track name and its
# Description text in row 6 of the Music Tracks table; the other columns
# (Major/Moderate/Minor Motifs) for that row are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = "The Mirthless"
$ws.Range("B6").Value = "Plays when the Mirthless fights the player. The Mirthless often haunts players even as they only begin to become insane, but at low sanity levels, it finally confronts the player."

# Matches the author's final cursor position recorded in the saved file.
$ws.Range("A7").Select()
